$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Table 2")
$ws3 = $wb.Worksheets.Item("Table 3")

$ws2.Range("B3").Value = "alpine skiing/snowboarding: 59% (n = 180)`nski touring/freeride: 3% (n = 9)`ncross-country skiing: 5.6% (n = 17)`nsledding: 4% (n = 12)`nice climbing: 0.33% (n = 1)`nhiking: 5.6% (n = 17)`nclimbing: 3.6% (n = 11)`nmountaineering: 0.66% (n = 2)`nbiking: 16% (n = 48)`nair sports: 0.33% (n = 1)`nwater sports: 0.33% (n = 1)`nother: 1.3% (n = 4)`nn = 303"

$ws3.Range("A8").Value = "Clinically relevant somatizaton symptoms (PHQ-15 ≥10)"
$ws3.Range("B8").Value = "5.9% (n = 18)"
